# Normalize multi-line shared-string cell text into single-line text
# (embedded line breaks -> single space), per the commit's cleanup pass.
# Also fixes a stray duplicate string on the "Adult Influenza Vaccine "
# sheet (B7) that should have referenced the same text already used
# elsewhere ("Fluvirin Preservative-free").

$wb = $excel.ActiveWorkbook

# --- Pediatric Vaccine sheet ---
$ws = $wb.Sheets.Item("Pediatric Vaccine ")
$ws.Range("D8").Value = "5 pack - 1 dose T-L syringes. No Needle"

# --- Adult Vaccine sheet ---
$ws = $wb.Sheets.Item("Adult Vaccine ")
$ws.Range("B14").Value = "Tetanus  Diphtheria Toxoids Adsorbed for Adults No Preservative"

# --- Pediatric Influenza Vaccine sheet ---
$ws = $wb.Sheets.Item("Pediatric Influenza Vaccine ")
$ws.Range("B3").Value = "Fluzone Pediatric dose No Preservative"
$ws.Range("B6").Value = "Fluarix Preservative-Free"
$ws.Range("B9").Value = "FluMist No Preservative"
$ws.Range("B10").Value = "Afluria No Preservative"
$ws.Range("H10").Value = "Merck (CSL product)"
$ws.Range("H11").Value = "Merck (CSL product)"
$ws.Range("B12").Value = "Afluria No Preservative"
$ws.Range("H12").Value = "Merck (CSL product)"

# --- Adult Influenza Vaccine sheet ---
$ws = $wb.Sheets.Item("Adult Influenza Vaccine ")
$ws.Range("B5").Value = "Agriflu No Preservative"
$ws.Range("B7").Value = "Fluvirin Preservative-free"
$ws.Range("B8").Value = "Fluarix Preservative-free"
$ws.Range("B10").Value = "Flumist No Preservative"
